$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.889.39'
$ws.Range("E2").Value = '  -2.25%  '

# Row 3
$ws.Range("D3").Value = '2.575.75'
$ws.Range("E3").Value = '  -4.07%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '549.15'
$ws.Range("E5").Value = '  -1.02%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.65'
$ws.Range("E6").Value = '  -2.35%  '

# Row 7
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.597'
$ws.Range("E8").Value = '  +2.63%  '

# Row 9
$ws.Range("E9").Value = '  -1.16%  '

# Row 10
$ws.Range("E10").Value = '  -1.39%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.43'
$ws.Range("E11").Value = '  +1.54%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.365'
$ws.Range("E12").Value = '  -0.75%  '

# Row 13
$ws.Range("D13").Value = '3.032.84'
$ws.Range("E13").Value = '  -4.13%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.58'
$ws.Range("E14").Value = '  -2.67%  '

# Row 15
$ws.Range("D15").Value = '61.782.31'
$ws.Range("E15").Value = '  -2.19%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000144'
$ws.Range("E16").Value = '  +0.08%  '

# Row 17
$ws.Range("D17").Value = '2.580.59'
$ws.Range("E17").Value = '  -4.15%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.53'
$ws.Range("E18").Value = '  -4.16%  '

# Row 19
$ws.Range("E19").Value = '  -0.43%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '337.48'
$ws.Range("E20").Value = '  -1.46%  '

# Row 21
$ws.Range("E21").Value = '  -4.21%  '

# Row 22
$ws.Range("E22").Value = '  +0.32%  '

# Row 23
$ws.Range("E23").Value = '  -2.44%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.39'
$ws.Range("E24").Value = '  +0.02%  '

# Row 25
$ws.Range("E25").Value = '  -0.43%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  +0.12%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.13'
$ws.Range("E27").Value = '  +1.23%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.45'
$ws.Range("E28").Value = '  +6.41%  '

# Row 29
$ws.Range("E29").Value = '  -1.66%  '

# Row 30
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.33'
$ws.Range("E30").Value = '  +0.19%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.90'
$ws.Range("E31").Value = '  -1.28%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '160.84'
$ws.Range("E32").Value = '  -2.66%  '

# Row 33
$ws.Range("E33").Value = '  +0.01%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.80'
$ws.Range("E34").Value = '  +1.44%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '19.20'
$ws.Range("E35").Value = '  -1.61%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.41'
$ws.Range("E36").Value = '  -0.92%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.79'
$ws.Range("E37").Value = '  +1.39%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '331.46'
$ws.Range("E38").Value = '  -2.36%  '

# Row 39
$ws.Range("E39").Value = '  -1.97%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.96'
$ws.Range("E40").Value = '  -1.87%  '

# Row 41
$ws.Range("E41").Value = '  +0.02%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '37.53'
$ws.Range("E42").Value = '  -1.41%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.77'
$ws.Range("E43").Value = '  +0.09%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("E44").Value = '  +0.00%  '

# Row 45
$ws.Range("D45").Value = '2.126.48'
$ws.Range("E45").Value = '  +1.36%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.606'
$ws.Range("E46").Value = '  -2.12%  '

# Row 47
$ws.Range("E47").Value = '  -1.07%  '

# Row 48
$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0546'
$ws.Range("E48").Value = '  -2.38%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.56'
$ws.Range("E49").Value = '  -3.27%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0964'
$ws.Range("E50").Value = '  -0.61%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0239'
$ws.Range("E51").Value = '  -0.84%  '
